$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update 2025-11 (row 24) stats
$ws.Range("B24").Value = 6430
$ws.Range("C24").Value = 1002
$ws.Range("D24").Value = 5987841
$ws.Range("E24").Value = 931.2349922239503
$ws.Range("F24").Value = 9.614728946471196
$ws.Range("G24").Value = 3.83419689119171
$ws.Range("H24").Value = 26.84566616079578
